# Updates crypto price/volume data per the commit:
# "Updated symbol list on Thu Feb  9 19:35:53 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '313.10'

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '-4.20%'

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '41.50'

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '-6.38%'

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.145'

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '-1.42%'

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.07920'

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '-5.41%'

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '4.360'

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '-1.94%'

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.648'

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '-14.90%'

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.9183'

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '-5.64%'

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1112'

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '-3.44%'

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.1810'

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '-5.27%'

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.09169'

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '-5.17%'

$ws.Range("B12").Value = 'BitrueCoin'

$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.04501'

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '-2.66%'

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '7.307'

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '-15.56%'

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.1052'

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '-0.69%'

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.001267'

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '-1.98%'

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.005948'

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '2.10%'

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.352'

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '-1.44%'

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '2.58%'

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.1392'

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '2.14%'

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '2.28%'

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.04160'

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '0.08%'

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.001250'

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '1.19%'

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.004191'

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '-5.39%'

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.0001229'

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '-5.76%'

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.0003005'

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '0.60%'

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02491'

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '-8.87%'

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.05313'

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '-5.79%'

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.008086'

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '2.69%'

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.1365'

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '-3.30%'

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.007648'

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '4.26%'

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.002076'

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '1.40%'

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.007541'

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '-4.94%'

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.3128'

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '-10.73%'

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00006803'

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '-0.79%'

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000757'

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '0.61%'

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.003409'

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '-2.46%'

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.004136'

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '16.85%'

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.00002119'

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '0.61%'

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0002018'

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '0.61%'
